$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input")

# Add elasticity (e_trend, column B) and shock (e_cycle, column C) to FDG accounts
# Row 2 = alcohol_taxes
$ws.Range("B2").Value = 0
# Row 3 = permits
$ws.Range("C3").Value = 1
# Row 4 = miscellaneous_income
$ws.Range("C4").Value = 1
# Row 5 = hydro_indexation
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1

# Move the active selection
$ws.Range("A4").Select()
